# Planilla de Trabajo HH - update hourly rate, add make-up hours, refresh selections.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HH")

# Updated "valor HH total" rate (was 217.38).
$ws.Range("S3").Formula = "=(284.03)*S2"

# New make-up-hour entries for Juan Carlos Garcés (mirrors R8's time format).
$ws.Range("R9").Value = 0.98472222222222217
$ws.Range("R10").Value = 0.84305555555555556
$ws.Range("R11").Value = 0.93958333333333333
$ws.Range("R9:R11").NumberFormat = $ws.Range("R8").NumberFormat

# Roll the new hours into the running total.
$ws.Range("Q6").Formula = "= Q8+R8+R9+R10+R11"

# Touched while reviewing row 20 - leaves a formatted-but-empty cell behind.
$ws.Range("R20").NumberFormat = $ws.Range("C29").NumberFormat

# Restore cursor positions: HH!S3 (even though Presupuesto stays the active tab)
# and Presupuesto!B11.
$pres = $wb.Worksheets.Item("Presupuesto")
[void]$pres.Range("B11").Select()
[void]$ws.Range("S3").Select()
[void]$pres.Activate()
